# Updated cryptos list on Fri Sep 22 09:56:09 UTC 2023 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) values for the
# coin rows on Sheet1, matching the latest scrape. Rows/columns that are
# unchanged in the source data are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value2 = '26.693.64'
$ws.Cells.Item(2, 5).Value2 = '  -0.40%  '

$ws.Cells.Item(3, 4).Value2 = '1.597.61'
$ws.Cells.Item(3, 5).Value2 = '  -0.98%  '

$ws.Cells.Item(4, 5).Value2 = '  -0.10%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value2 = '211.47'
$ws.Cells.Item(5, 5).Value2 = '  -0.25%  '

$ws.Cells.Item(7, 5).Value2 = '  -0.19%  '

$ws.Cells.Item(8, 5).Value2 = '  -0.60%  '

$ws.Cells.Item(9, 5).Value2 = '  -1.36%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value2 = '19.75'
$ws.Cells.Item(10, 5).Value2 = '  +0.40%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value2 = '0.0839'
$ws.Cells.Item(11, 5).Value2 = '  +0.23%  '

$ws.Cells.Item(12, 4).Value2 = '1.821.31'
$ws.Cells.Item(12, 5).Value2 = '  -1.07%  '

$ws.Cells.Item(13, 4).Value2 = '1.594.87'
$ws.Cells.Item(13, 5).Value2 = '  -1.29%  '

$ws.Cells.Item(14, 5).Value2 = '  -1.07%  '

$ws.Cells.Item(15, 5).Value2 = '  -1.84%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value2 = '65.09'
$ws.Cells.Item(16, 5).Value2 = '  +1.91%  '

$ws.Cells.Item(17, 4).Value2 = '26.699.06'
$ws.Cells.Item(17, 5).Value2 = '  -0.47%  '

$ws.Cells.Item(18, 5).Value2 = '  -0.38%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value2 = '210.09'
$ws.Cells.Item(19, 5).Value2 = '  -0.17%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value2 = '6.75'
$ws.Cells.Item(21, 5).Value2 = '  -0.01%  '

$ws.Cells.Item(22, 5).Value2 = '  -0.58%  '

$ws.Cells.Item(23, 5).Value2 = '  -1.21%  '

$ws.Cells.Item(24, 5).Value2 = '  -0.12%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value2 = '146.64'
$ws.Cells.Item(25, 5).Value2 = '  +0.15%  '

$ws.Cells.Item(26, 5).Value2 = '  -0.10%  '

$ws.Cells.Item(27, 5).Value2 = '  -4.11%  '

$ws.Cells.Item(28, 5).Value2 = '  +2.03%  '

$ws.Cells.Item(29, 5).Value2 = '  -0.42%  '

$ws.Cells.Item(30, 5).Value2 = '  +0.13%  '

$ws.Cells.Item(31, 5).Value2 = '  -0.70%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value2 = '3.21'
$ws.Cells.Item(32, 5).Value2 = '  -1.35%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value2 = '0.663'
$ws.Cells.Item(33, 5).Value2 = '  -6.22%  '

$ws.Cells.Item(34, 5).Value2 = '  -1.49%  '

$ws.Cells.Item(35, 4).Value2 = '1.298.34'
$ws.Cells.Item(35, 5).Value2 = '  -1.83%  '

$ws.Cells.Item(36, 5).Value2 = '  -0.17%  '

$ws.Cells.Item(37, 5).Value2 = '  -3.81%  '

$ws.Cells.Item(38, 5).Value2 = '  -1.12%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value2 = '0.843'
$ws.Cells.Item(39, 5).Value2 = '  +1.91%  '

$ws.Cells.Item(40, 5).Value2 = '  -0.10%  '

$ws.Cells.Item(41, 5).Value2 = '  +1.54%  '

$ws.Cells.Item(42, 5).Value2 = '  -0.49%  '

$ws.Cells.Item(43, 5).Value2 = '  -0.44%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value2 = '63.74'
$ws.Cells.Item(44, 5).Value2 = '  +0.46%  '

$ws.Cells.Item(45, 4).Value2 = '1.733.90'
$ws.Cells.Item(45, 5).Value2 = '  -1.10%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value2 = '0.888'
$ws.Cells.Item(46, 5).Value2 = '  +8.99%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value2 = '90.18'
$ws.Cells.Item(47, 5).Value2 = '  +1.12%  '

$ws.Cells.Item(48, 5).Value2 = '  +0.83%  '

$ws.Cells.Item(49, 5).Value2 = '  +2.38%  '

$ws.Cells.Item(50, 5).Value2 = '  -1.47%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value2 = '7.47'
$ws.Cells.Item(51, 5).Value2 = '  -0.22%  '
